$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row to match the new API field names
$ws.Range("A1").Value = "full_name"
$ws.Range("B1").Value = "cnic"
$ws.Range("C1").Value = "check_in"
$ws.Range("D1").Value = "check_out"
$ws.Range("E1").Value = "user_id"

# Append the new user record created through the "create user" API
$ws.Range("A6").Value = "aahil alwani"
$ws.Range("B6").Value = "42w322e633333"
$ws.Range("C6").Value = "30 apr 2025 1:02pm |"

# user_id is stored as text (it is a shared-string cell in the source data),
# so force text formatting before writing the numeric-looking id, then drop
# back to the default style so no stray formatting is left on the cell.
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1232536327"
$ws.Range("E6").Style = "Normal"

# Match the new selection left behind by the editing session
$ws.Range("A10").Select()

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
